$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the date in A1 (serial 45406 -> 45436, i.e. 2024-04-24 -> 2024-05-24)
$ws.Range("A1").Value = 45436

# Update the price values in column D
$ws.Range("D33").Value = 1507.712
$ws.Range("D34").Value = 2542.857
$ws.Range("D35").Value = 2799.394
